# Update cryptos list: Price (D) and Volume(1h) (E) columns
# Leading apostrophe forces text type so numeric-looking prices
# (e.g. "1.000", "0.9992") are not coerced into Excel numbers,
# matching the original inlineStr text-cell semantics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.777.71"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "'1.875.73"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'0.7302"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "'241.62"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").Value = "'0.07111"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  -0.99%  "

$ws.Range("D11").Value = "'0.08257"
$ws.Range("E11").Value = "  -2.73%  "

$ws.Range("D13").Value = "'1.876.99"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "'5.324"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").Value = "'92.55"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").Value = "'29.809.84"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "'6.045"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "'248.05"
$ws.Range("E18").Value = "  +2.16%  "

$ws.Range("D19").Value = "'13.39"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").Value = "'0.000007824"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").Value = "'2.154.67"
$ws.Range("E21").Value = "  +1.89%  "

$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'0.9990"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "'7.738"
$ws.Range("E24").Value = "  -2.31%  "

$ws.Range("D25").Value = "'0.1542"
$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("D26").Value = "'9.178"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").Value = "'162.92"
$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "'1.438"
$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("D31").Value = "'4.540"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "'4.207"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("D34").Value = "'0.05277"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'1.236"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").Value = "'2.705"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("D40").Value = "'2.747"
$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("D41").Value = "'0.4492"
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("D42").Value = "'6.012"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").Value = "'0.8669"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'71.37"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").Value = "'1.064.18"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").Value = "'104.67"
$ws.Range("E46").Value = "  +2.03%  "

$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").Value = "'1.830"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").Value = "'7.501"
$ws.Range("E49").Value = "  -2.86%  "

$ws.Range("D50").Value = "'9.517"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("D51").Value = "'2.036.42"
$ws.Range("E51").Value = "  +0.77%  "
